$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15; this shifts existing rows 15..120 down to 16..121.
$ws.Rows(15).Insert()

# Populate the newly inserted row 15 with the new weekly record. The
# "constant" columns (A,B,C,E,F,G,H,I,N,O,Q,R) hold the same value on every
# data row of this sheet, so they are written back literally here.
$ws.Range("A15").Value = 8
$ws.Range("B15").Value = "Terminal La Palmera de La Serena"
$ws.Range("C15").Value = "Coquimbo"
$ws.Range("D15").Value = 44602
$ws.Range("E15").Value = 4
$ws.Range("F15").Value = 100112040
$ws.Range("G15").Value = "Cilantro"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 2200
$ws.Range("K15").Value = 2300
$ws.Range("L15").Value = 2500
$ws.Range("M15").Value = 2400
$ws.Range("N15").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O15").Value = "Provincia del Elquí"
$ws.Range("P15").Value = 1600
$ws.Range("Q15").Value = 1.5
$ws.Range("R15").Value = "Hortaliza"
